# Apply odds updates to Sheet1 (changed cells per the commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("U3").Value = 1.52
$ws.Range("V3").Value = 2.46

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.65
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 2.22
$ws.Range("K4").Value = 2.22
$ws.Range("L4").Value = 4.35
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 2.02
$ws.Range("U4").Value = 1.64
$ws.Range("V4").Value = 2.21
$ws.Range("W4").Value = 6.9
$ws.Range("X4").Value = 7.4
$ws.Range("Y4").Value = 6.9
$ws.Range("Z4").Value = 11.25
$ws.Range("AA4").Value = 10.5
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 11.5
$ws.Range("AF4").Value = 40
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 11.5
$ws.Range("AL4").Value = 28
$ws.Range("AM4").Value = 28
$ws.Range("AN4").Value = 3.65
$ws.Range("AO4").Value = 8.25
$ws.Range("AP4").Value = 16
$ws.Range("AQ4").Value = 27
$ws.Range("AR4").Value = 50
$ws.Range("AS4").Value = 175
$ws.Range("AT4").Value = 3
$ws.Range("AU4").Value = 7
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 6.1
$ws.Range("AX4").Value = 22
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 120
$ws.Range("BA4").Value = 150
$ws.Range("BB4").Value = 300

# Row 5
$ws.Range("G5").Value = 2.38
$ws.Range("I5").Value = 2.6
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 3.25
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 15
$ws.Range("X5").Value = 13
$ws.Range("Z5").Value = 23
$ws.Range("AB5").Value = 23
$ws.Range("AG5").Value = 126
$ws.Range("AK5").Value = 29
$ws.Range("AL5").Value = 21
$ws.Range("AW5").Value = 5
$ws.Range("AX5").Value = 15

# Row 9
$ws.Range("G9").Value = 2.05
$ws.Range("I9").Value = 3.8
$ws.Range("X9").Value = 9.5
$ws.Range("AL9").Value = 29
$ws.Range("BB9").Value = 201

# Row 10
$ws.Range("Z10").Value = 15
$ws.Range("AK10").Value = 41
$ws.Range("BB10").Value = 151

# Row 11
$ws.Range("G11").Value = 4.5
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 1.7
$ws.Range("J11").Value = 5
$ws.Range("K11").Value = 2.12
$ws.Range("L11").Value = 2.3
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 6.6
$ws.Range("O11").Value = 1.37
$ws.Range("P11").Value = 2.85
$ws.Range("Q11").Value = 2.1
$ws.Range("R11").Value = 1.65
$ws.Range("T11").Value = 2.65
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.72
$ws.Range("W11").Value = 10.75
$ws.Range("X11").Value = 25
$ws.Range("Y11").Value = 15.5
$ws.Range("AA11").Value = 50
$ws.Range("AB11").Value = 60
$ws.Range("AC11").Value = 6.6
$ws.Range("AD11").Value = 6.8
$ws.Range("AE11").Value = 18
$ws.Range("AF11").Value = 100
$ws.Range("AH11").Value = 5.9
$ws.Range("AI11").Value = 7.3
$ws.Range("AJ11").Value = 8.5
$ws.Range("AK11").Value = 13
$ws.Range("AL11").Value = 15
$ws.Range("AM11").Value = 32
$ws.Range("AR11").Value = 250
$ws.Range("AT11").Value = 2.65
$ws.Range("AY11").Value = 19.5
